$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $value) {
    # Force the value to be stored as text, even if it looks like a number,
    # matching the source workbook where these are text labels (e.g. "1.00").
    $r = $ws.Range($cellRange)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

$ws.Range("D2").Value = "60.675.42"
$ws.Range("E2").Value = "  +3.27%  "

$ws.Range("D3").Value = "2.686.41"
$ws.Range("E3").Value = "  +1.38%  "

Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.25%  "

Set-TextValue "D5" "523.63"
$ws.Range("E5").Value = "  +2.00%  "

Set-TextValue "D6" "146.47"
$ws.Range("E6").Value = "  +1.85%  "

$ws.Range("E7").Value = "  +0.13%  "

Set-TextValue "D8" "0.576"
$ws.Range("E8").Value = "  +1.68%  "

$ws.Range("D9").Value = "2.705.04"
$ws.Range("E9").Value = "  +0.99%  "

Set-TextValue "D10" "6.44"
$ws.Range("E10").Value = "  +2.69%  "

$ws.Range("E11").Value = "  +0.29%  "

$ws.Range("E12").Value = "  +2.07%  "

$ws.Range("E13").Value = "  +1.81%  "

$ws.Range("D14").Value = "3.158.08"
$ws.Range("E14").Value = "  +1.10%  "

$ws.Range("D15").Value = "60.434.01"
$ws.Range("E15").Value = "  +2.81%  "

Set-TextValue "D16" "21.31"
$ws.Range("E16").Value = "  +1.26%  "

$ws.Range("D17").Value = "2.779.79"
$ws.Range("E17").Value = "  +3.64%  "

$ws.Range("E18").Value = "  +1.77%  "

Set-TextValue "D19" "351.23"
$ws.Range("E19").Value = "  +2.67%  "

$ws.Range("E20").Value = "  +0.26%  "

Set-TextValue "D21" "10.63"
$ws.Range("E21").Value = "  +2.00%  "

$ws.Range("E22").Value = "  +3.15%  "

$ws.Range("E23").Value = "  +0.24%  "

Set-TextValue "D24" "62.82"
$ws.Range("E24").Value = "  +2.96%  "

Set-TextValue "D25" "0.422"
$ws.Range("E25").Value = "  +0.69%  "

$ws.Range("E26").Value = "  +5.46%  "

Set-TextValue "D27" "0.995"
$ws.Range("E27").Value = "  +0.39%  "

$ws.Range("D28").Value = "0.0₃0818"
$ws.Range("E28").Value = "  +1.19%  "

$ws.Range("E29").Value = "  +0.97%  "

$ws.Range("E30").Value = "  +7.45%  "

$ws.Range("E31").Value = "  +0.04%  "

$ws.Range("E32").Value = "  +1.48%  "

Set-TextValue "D33" "19.10"
$ws.Range("E33").Value = "  +1.14%  "

Set-TextValue "D34" "148.27"
$ws.Range("E34").Value = "  -0.81%  "

$ws.Range("E35").Value = "  +7.62%  "

$ws.Range("E36").Value = "  +8.24%  "

Set-TextValue "D37" "0.951"
$ws.Range("E37").Value = "  -6.04%  "

$ws.Range("E38").Value = "  +10.82%  "

Set-TextValue "D39" "0.879"
$ws.Range("E39").Value = "  +3.19%  "

Set-TextValue "D40" "36.92"
$ws.Range("E40").Value = "  +1.15%  "

Set-TextValue "D41" "3.69"
$ws.Range("E41").Value = "  +0.69%  "

Set-TextValue "D42" "282.06"
$ws.Range("E42").Value = "  +0.30%  "

$ws.Range("E43").Value = "  -0.94%  "

Set-TextValue "D44" "0.997"
$ws.Range("E44").Value = "  +0.38%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D45" "0.0988"
$ws.Range("E45").Value = "  +1.02%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D46" "19.95"
$ws.Range("E46").Value = "  +1.81%  "

$ws.Range("D47").Value = "2.135.52"
$ws.Range("E47").Value = "  +7.06%  "

$ws.Range("E48").Value = "  +1.30%  "

Set-TextValue "D49" "4.88"
$ws.Range("E49").Value = "  +3.91%  "

$ws.Range("E50").Value = "  +2.44%  "

Set-TextValue "D51" "10.46"
$ws.Range("E51").Value = "  +1.75%  "

